# Removing less than USD 5 price from extrapolation calibration because it
# is just a noise. This updates the recalculated calibration outputs
# (ABSM1_RN, M1_RN, CM2_RN, CMN3_RN, CMN4_RN) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row -> (D, E, F, G, H)
$updates = @{
    6  = @(116311.8085317966, -0.06347546924066284, 0.2630818951983989, -1.148958191313806, 9.973416671127778)
    8  = @(117984.6608108865, -0.07060723159594018, 0.2300712166986094, -0.9941041241852699, 7.116001874476617)
    9  = @(119745.3001773226, -0.09930866618743497, 0.3712322134719819, -2.078797047244133, 13.11549548757357)
    10 = @(121209.9761724819, -0.1334821789809024, 0.455209889918833, -1.923092694554094, 9.561231332881119)
    11 = @(123185.7306783697, -0.2117033195092484, 0.7983813384258449, -2.569052455221928, 12.27765383625103)
    12 = @(113569.9061873077, -0.1757385252146718, 0.1635854293001148, -1.388265908861001, 9.496335272112061)
    13 = @(113608.3707815903, -0.1029819612380781, 0.1548964532373811, -0.8248907173442492, 6.782178154826677)
    14 = @(113556.7920465209, -0.1555273097598965, 0.166218193832818, -1.745819275818429, 13.76211482537033)
    15 = @(113619.7383582534, -0.1359635606531136, 0.1586718942005126, -1.015307318833194, 8.886208270400793)
    17 = @(113554.6811280577, -0.1208637641525671, 0.1488618419469296, -1.148920369578007, 9.843347091434062)
    18 = @(113668.5360785112, -0.08675937881418373, 0.1394736117780836, -0.5656137110319629, 6.227484561738295)
    20 = @(114315.533801794, -0.05863437011426385, 0.1999913316923044, -1.553362806442369, 13.16252457306894)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
    $ws.Cells.Item($row, 7).Value = $vals[3]
    $ws.Cells.Item($row, 8).Value = $vals[4]
}
